$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new F-column (dSF) value, reflecting a repull/recalculation of data.
$updates = @{
    2  = -1
    7  = -2
    12 = 2
    22 = -2
    30 = -7
    31 = 0
    34 = 1
    35 = 4
    37 = -2
    41 = 1
    45 = 0
    47 = 3
    49 = 3
    51 = 1
    53 = 2
    54 = 3
    55 = 3
    56 = 1
    63 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
